# Fixed update to excel issue
#
# 1. Rename "Requested quantity" header to "Weekly_PO_Qty" on the
#    "Weekly Quantity" sheet and to "Monthly_PO_Qty" on the
#    "Monthly Trend" sheet.
# 2. Add a new "PO Forecast" sheet (after "Monthly Trend") containing
#    the ds / PO_Forecast / yhat_lower / yhat_upper forecast table.

$wb = $excel.ActiveWorkbook

$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# Add the new sheet right after "Monthly Trend" so ordering matches
# Weekly Quantity, Monthly Trend, PO Forecast.
$wsForecast = $wb.Worksheets.Add($null, $wsMonthly)
$wsForecast.Name = "PO Forecast"

$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"
$wsForecast.Range("A1:D1").Font.Bold = $true

$forecastData = @(
    @(44948.99999999999, 67, 52.27319200548753, 81.70606421874353),
    @(44969.99999999999, 36, 20.82558191678516, 53.27988143470169),
    @(44976.99999999999, 26, 10.634085979578, 41.38664876761712),
    @(44983.99999999999, 15, -0.7029172815474665, 31.81217290933887),
    @(44990.99999999999, 5, -10.08045653808977, 21.69765041515405),
    @(44997.99999999999, 0, -21.74444416223787, 9.662522766955924),
    @(45004.99999999999, 0, -31.59819658987964, 0.9049957939959034),
    @(45011.99999999999, 0, -41.60279389909932, -10.13022346088148),
    @(45018.99999999999, 0, -52.42109792703413, -19.03760917540954),
    @(45025.99999999999, 0, -61.48257682409626, -30.6355719172887),
    @(45032.99999999999, 0, -72.3361733210949, -39.96107352396491),
    @(45039.99999999999, 0, -82.22584766320907, -50.18875579556106)
)

$row = 2
foreach ($rec in $forecastData) {
    $wsForecast.Cells.Item($row, 1).Value = $rec[0]
    $wsForecast.Cells.Item($row, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $wsForecast.Cells.Item($row, 2).Value = $rec[1]
    $wsForecast.Cells.Item($row, 3).Value = $rec[2]
    $wsForecast.Cells.Item($row, 4).Value = $rec[3]
    $row++
}

$wsForecast.Range("A1").Select() | Out-Null
